$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 5620
$ws.Range("I69").Value = 5740
$ws.Range("J69").Value = 5524
$ws.Range("K69").Value = 17220
$ws.Range("L69").Value = 16572
$ws.Range("M69").Value = -16346
$ws.Range("N69").Value = -18320
$ws.Range("H72").Value = 5620
$ws.Range("I72").Value = 5740
$ws.Range("J72").Value = 5524
$ws.Range("K72").Value = 51660
$ws.Range("L72").Value = 49716
$ws.Range("M72").Value = -47292
$ws.Range("N72").Value = -58452
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H113").Value = 3157.9167
$ws.Range("I113").Value = 2505.8333
$ws.Range("J113").Value = 3810
$ws.Range("K113").Value = 2505.8333
$ws.Range("L113").Value = 3810
$ws.Range("M113").Value = 748.1667000000002
$ws.Range("N113").Value = -10318
$ws.Range("H129").Value = 1478.0638
$ws.Range("I129").Value = 645.7857
$ws.Range("J129").Value = 1831.1515
$ws.Range("K129").Value = 1937.3571
$ws.Range("L129").Value = 5493.4545
$ws.Range("M129").Value = 3062.6429
$ws.Range("N129").Value = -15493.4545

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2664.69
$ws.Range("I32").Value = 2664.69
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2664.69
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2377.69
$ws.Range("N32").ClearContents()
$ws.Range("H57").Value = 50000
$ws.Range("I57").Value = 50000
$ws.Range("K57").Value = 50000
$ws.Range("M57").Value = -49516
$ws.Range("H74").Value = 3559.36
$ws.Range("I74").Value = 4672.2666
$ws.Range("J74").Value = 1890
$ws.Range("K74").Value = 4672.2666
$ws.Range("L74").Value = 1890
$ws.Range("M74").Value = -3798.2666
$ws.Range("N74").Value = -3638
$ws.Range("H77").Value = 3559.36
$ws.Range("I77").Value = 4672.2666
$ws.Range("J77").Value = 1890
$ws.Range("K77").Value = 23361.333
$ws.Range("L77").Value = 9450
$ws.Range("M77").Value = -18993.333
$ws.Range("N77").Value = -18186
$ws.Range("H102").Value = 2645.1
$ws.Range("I102").Value = 1374.2
$ws.Range("J102").Value = 3916
$ws.Range("K102").Value = 1374.2
$ws.Range("L102").Value = 3916
$ws.Range("M102").Value = 247.8
$ws.Range("N102").Value = -7160
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H132").Value = 3271.9624
$ws.Range("I132").Value = 1886.4166
$ws.Range("J132").Value = 5350.2812
$ws.Range("K132").Value = 5659.2498
$ws.Range("L132").Value = 16050.8436
$ws.Range("M132").Value = -3129.2498
$ws.Range("N132").Value = -21110.8436

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 53750
$ws.Range("J63").Value = 53750
$ws.Range("L63").Value = 53750
$ws.Range("N63").Value = -55122
$ws.Range("H66").Value = 53750
$ws.Range("J66").Value = 53750
$ws.Range("L66").Value = 161250
$ws.Range("N66").Value = -168114
$ws.Range("H99").Value = 2275.6667
$ws.Range("I99").Value = 1921
$ws.Range("J99").Value = 2985
$ws.Range("K99").Value = 1921
$ws.Range("L99").Value = 2985
$ws.Range("M99").Value = -423
$ws.Range("N99").Value = -5981
$ws.Range("H103").Value = 12789.25
$ws.Range("J103").Value = 12789.25
$ws.Range("L103").Value = 12789.25
$ws.Range("N103").Value = -15133.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 18307.2
$ws.Range("J28").Value = 18307.2
$ws.Range("L28").Value = 18307.2
$ws.Range("N28").Value = -18797.2
$ws.Range("H31").Value = 1802.58
$ws.Range("I31").Value = 1059.9844
$ws.Range("J31").Value = 3122.75
$ws.Range("K31").Value = 1059.9844
$ws.Range("L31").Value = 3122.75
$ws.Range("M31").Value = -764.9844000000001
$ws.Range("N31").Value = -3712.75
$ws.Range("H34").Value = 1802.58
$ws.Range("I34").Value = 1059.9844
$ws.Range("J34").Value = 3122.75
$ws.Range("K34").Value = 1059.9844
$ws.Range("L34").Value = 3122.75
$ws.Range("M34").Value = -857.9844000000001
$ws.Range("N34").Value = -3526.75
$ws.Range("H43").Value = 20100
$ws.Range("J43").Value = 20100
$ws.Range("L43").Value = 20100
$ws.Range("N43").Value = -20468
$ws.Range("H58").Value = 6270.2085
$ws.Range("I58").Value = 3078.2144
$ws.Range("J58").Value = 10739
$ws.Range("K58").Value = 3078.2144
$ws.Range("L58").Value = 10739
$ws.Range("M58").Value = -2875.2144
$ws.Range("N58").Value = -11145
$ws.Range("H74").Value = 18616.857
$ws.Range("J74").Value = 18616.857
$ws.Range("L74").Value = 18616.857
$ws.Range("N74").Value = -20364.857
$ws.Range("H77").Value = 18616.857
$ws.Range("J77").Value = 18616.857
$ws.Range("L77").Value = 55850.571
$ws.Range("N77").Value = -64586.571
$ws.Range("H101").Value = 20100
$ws.Range("J101").Value = 20100
$ws.Range("L101").Value = 20100
$ws.Range("N101").Value = -26590
$ws.Range("H132").Value = 2481.56
$ws.Range("I132").Value = 1667.9615
$ws.Range("J132").Value = 3362.9583
$ws.Range("K132").Value = 5003.8845
$ws.Range("L132").Value = 10088.8749
$ws.Range("M132").Value = -2473.8845
$ws.Range("N132").Value = -15148.8749
$ws.Range("H134").Value = 2696.6
$ws.Range("I134").Value = 1773.8
$ws.Range("J134").Value = 3158
$ws.Range("K134").Value = 5321.4
$ws.Range("L134").Value = 9474
$ws.Range("M134").Value = -2786.4
$ws.Range("N134").Value = -14544
$ws.Range("H136").Value = 6270.2085
$ws.Range("I136").Value = 3078.2144
$ws.Range("J136").Value = 10739
$ws.Range("K136").Value = 9234.643199999999
$ws.Range("L136").Value = 32217
$ws.Range("M136").Value = -6684.643199999999
$ws.Range("N136").Value = -37317

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 3463.3333
$ws.Range("I9").Value = 2450
$ws.Range("J9").Value = 5490
$ws.Range("K9").Value = 7350
$ws.Range("L9").Value = 16470
$ws.Range("M9").Value = -7126
$ws.Range("N9").Value = -16918

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2408.7222
$ws.Range("I80").Value = 2375.3333
$ws.Range("J80").Value = 2475.5
$ws.Range("K80").Value = 2375.3333
$ws.Range("L80").Value = 2475.5
$ws.Range("M80").Value = -1377.3333
$ws.Range("N80").Value = -4471.5
$ws.Range("H83").Value = 2408.7222
$ws.Range("I83").Value = 2375.3333
$ws.Range("J83").Value = 2475.5
$ws.Range("K83").Value = 11876.6665
$ws.Range("L83").Value = 12377.5
$ws.Range("M83").Value = -6884.666499999999
$ws.Range("N83").Value = -22361.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2135.2703
$ws.Range("I16").Value = 1252.5385
$ws.Range("K16").Value = 1252.5385
$ws.Range("M16").Value = -1082.5385
$ws.Range("H68").Value = 1905.4166
$ws.Range("I68").Value = 1586.8422
$ws.Range("J68").Value = 3116
$ws.Range("K68").Value = 1586.8422
$ws.Range("L68").Value = 3116
$ws.Range("M68").Value = -837.8422
$ws.Range("N68").Value = -4614
$ws.Range("H71").Value = 1905.4166
$ws.Range("I71").Value = 1586.8422
$ws.Range("J71").Value = 3116
$ws.Range("K71").Value = 7934.211
$ws.Range("L71").Value = 15580
$ws.Range("M71").Value = -4190.211
$ws.Range("N71").Value = -23068

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 5626
$ws.Range("J74").Value = 5626
$ws.Range("L74").Value = 5626
$ws.Range("N74").Value = -7498
$ws.Range("H77").Value = 5626
$ws.Range("J77").Value = 5626
$ws.Range("L77").Value = 16878
$ws.Range("N77").Value = -26238
$ws.Range("H100").Value = 569.35
$ws.Range("I100").Value = 453
$ws.Range("J100").Value = 840.8333
$ws.Range("K100").Value = 906
$ws.Range("L100").Value = 1681.6666
$ws.Range("M100").Value = -365
$ws.Range("N100").Value = -2763.6666
$ws.Range("H132").Value = 1477.2678
$ws.Range("I132").Value = 1184.0238
$ws.Range("K132").Value = 3552.0714
$ws.Range("M132").Value = -1022.0714
$ws.Range("H136").Value = 2132.3489
$ws.Range("I136").Value = 1609.2258
$ws.Range("J136").Value = 3483.75
$ws.Range("K136").Value = 4827.6774
$ws.Range("L136").Value = 10451.25
$ws.Range("M136").Value = -2277.6774
$ws.Range("N136").Value = -15551.25
